# Add (Ito & Claesson-Welsh, 1999) to VEGF:VEGFR1
# The VEGFA165_VEGFR1 sheet's table previously ended at row 9 (Breier et al., 1995).
# We insert a new last row holding the Ito & Claesson-Welsh, 1999 datapoint, which
# means row 9's "closing" (thick double-bottom-border) formatting needs to move down
# to the new row 10, while row 9 itself becomes a normal interior row (same look as
# rows 7-8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEGFA165_VEGFR1")

# 1) Push the current last row's formatting down onto the brand-new row 10 first
#    (copies both the values and the "table bottom" border/height formatting that
#    currently lives on row 9).
$ws.Range("A9:D9").Copy($ws.Range("A10:D10"))

# 2) Overwrite row 10 with the new reference's actual data.
$ws.Range("A10").Value = "Ito & Claesson-Welsh, 1999"
$ws.Range("B10").Value = "Radioligand"
$ws.Range("C10").Value = 91.54
$ws.Range("D10").Value = ""

# 3) Re-style row 9 as a normal interior row by copying the format from row 8
#    (the row right above it, which already has the interior-row look).
$ws.Range("A8:D8").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(9).AutoFit()              # drop the now-stale thick-bottom row height
$excel.CutCopyMode = $false

# 4) Make the VEGFA165_VEGFR1 tab the active sheet/selection, matching the saved view.
$ws.Activate() | Out-Null
$ws.Range("B14").Select() | Out-Null
